$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (avoid Excel auto-converting numeric-looking /
# date-looking strings into numbers / dates) for the cells being rewritten.
$ws.Range("A2:G4").NumberFormat = "@"
$ws.Range("K2:L4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "230"
$ws.Range("B2").Value = "220.8 - 239.2"
$ws.Range("C2").Value = "227.8"
$ws.Range("D2").Value = "231.4"
$ws.Range("E2").Value = "234.2"
$ws.Range("F2").Value = "234.6"
$ws.Range("G2").Value = "233.0"
$ws.Range("K2").Value = "2025-03-05"
$ws.Range("L2").Value = "2026-03-05"

# Row 3
$ws.Range("A3").Value = "150"
$ws.Range("B3").Value = "144.0 - 156.0"
$ws.Range("C3").Value = "145.7"
$ws.Range("D3").Value = "151.5"
$ws.Range("E3").Value = "154.8"
$ws.Range("F3").Value = "148.6"
$ws.Range("G3").Value = "149.7"
$ws.Range("K3").Value = "2025-03-05"
$ws.Range("L3").Value = "2026-03-05"

# Row 4
$ws.Range("A4").Value = "70"
$ws.Range("B4").Value = "67.2 - 72.8"
$ws.Range("C4").Value = "69.3"
$ws.Range("D4").Value = "69.9"
$ws.Range("E4").Value = "69.8"
$ws.Range("F4").Value = "69.9"
$ws.Range("G4").Value = "70.2"
$ws.Range("K4").Value = "2025-03-05"
$ws.Range("L4").Value = "2026-03-05"

# Remove column Q entirely (the "OCR Text" header and its now-unused data
# cells) and shrink the sheet's used range back down to A1:P4.
$ws.Range("Q1:Q4").EntireColumn.Delete()
